$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix separators in a few "Razon social" entries: commas -> periods ---
$ws.Range("E66").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E205").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E125").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("E135").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E148").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- Fix floating point numbers stored as text in the "Importe" column (H) ---
# Pre-format the whole range as Text so Excel keeps these as strings instead of
# re-parsing them into numbers (which would lose formatting / become General numbers).
$importeRange = $ws.Range("H2:H209")
$importeRange.NumberFormat = "@"

$ws.Range("H2").Value = "13948.00"
$ws.Range("H3").Value = "1480800.00"
$ws.Range("H4").Value = "880704.00"
$ws.Range("H5").Value = "2710.95"
$ws.Range("H6").Value = "496119.46"
$ws.Range("H7").Value = "166760.00"
$ws.Range("H8").Value = "343335.24"
$ws.Range("H9").Value = "985633.00"
$ws.Range("H10").Value = "304728.70"
$ws.Range("H11").Value = "56005.00"
$ws.Range("H12").Value = "1020566.00"
$ws.Range("H13").Value = "277843.74"
$ws.Range("H14").Value = "168393.39"
$ws.Range("H15").Value = "52000.00"
$ws.Range("H16").Value = "14940.00"
$ws.Range("H17").Value = "5060.00"
$ws.Range("H18").Value = "79457.70"
$ws.Range("H19").Value = "75.98"
$ws.Range("H20").Value = "2430329.73"
$ws.Range("H21").Value = "4400.00"
$ws.Range("H22").Value = "318628.21"
$ws.Range("H23").Value = "40998.00"
$ws.Range("H24").Value = "40784.12"
$ws.Range("H25").Value = "104470.00"
$ws.Range("H26").Value = "15493.00"
$ws.Range("H27").Value = "288852.20"
$ws.Range("H28").Value = "3663.49"
$ws.Range("H29").Value = "5500.00"
$ws.Range("H30").Value = "500.00"
$ws.Range("H31").Value = "10135.77"
$ws.Range("H32").Value = "2835.30"
$ws.Range("H33").Value = "8899.40"
$ws.Range("H34").Value = "1717.32"
$ws.Range("H35").Value = "16000.00"
$ws.Range("H36").Value = "18000.00"
$ws.Range("H37").Value = "9784.78"
$ws.Range("H38").Value = "12819.94"
$ws.Range("H39").Value = "21881.74"
$ws.Range("H40").Value = "6826.00"
$ws.Range("H41").Value = "2226.10"
$ws.Range("H42").Value = "1172.00"
$ws.Range("H43").Value = "20560.00"
$ws.Range("H44").Value = "1395.00"
$ws.Range("H45").Value = "4697.00"
$ws.Range("H46").Value = "5742.65"
$ws.Range("H47").Value = "1493.86"
$ws.Range("H48").Value = "1935.98"
$ws.Range("H49").Value = "2386.00"
$ws.Range("H50").Value = "244.40"
$ws.Range("H51").Value = "11048.00"
$ws.Range("H52").Value = "13100.00"
$ws.Range("H53").Value = "14556.47"
$ws.Range("H54").Value = "16000.00"
$ws.Range("H55").Value = "35000.00"
$ws.Range("H56").Value = "10800.00"
$ws.Range("H57").Value = "1935.52"
$ws.Range("H58").Value = "4200.00"
$ws.Range("H59").Value = "3300.00"
$ws.Range("H60").Value = "4700.00"
$ws.Range("H61").Value = "194567.00"
$ws.Range("H62").Value = "165180.00"
$ws.Range("H63").Value = "1240.00"
$ws.Range("H64").Value = "9500.00"
$ws.Range("H65").Value = "7980.00"
$ws.Range("H66").Value = "11610.00"
$ws.Range("H67").Value = "3910.00"
$ws.Range("H68").Value = "206456.25"
$ws.Range("H69").Value = "10920.00"
$ws.Range("H70").Value = "15960.42"
$ws.Range("H71").Value = "4600.00"
$ws.Range("H72").Value = "60380.00"
$ws.Range("H73").Value = "24000.00"
$ws.Range("H74").Value = "206.50"
$ws.Range("H75").Value = "15000.00"
$ws.Range("H76").Value = "34275.00"
$ws.Range("H77").Value = "71632.00"
$ws.Range("H78").Value = "2500.00"
$ws.Range("H79").Value = "9660.00"
$ws.Range("H80").Value = "1030.00"
$ws.Range("H81").Value = "6826.56"
$ws.Range("H82").Value = "18900.00"
$ws.Range("H83").Value = "6673.00"
$ws.Range("H84").Value = "73936.81"
$ws.Range("H85").Value = "3018.00"
$ws.Range("H86").Value = "33440.00"
$ws.Range("H87").Value = "840.00"
$ws.Range("H88").Value = "9655.00"
$ws.Range("H89").Value = "126.08"
$ws.Range("H90").Value = "19600.00"
$ws.Range("H91").Value = "6500.00"
$ws.Range("H92").Value = "5500.00"
$ws.Range("H93").Value = "9000.00"
$ws.Range("H94").Value = "2522.36"
$ws.Range("H95").Value = "44097.00"
$ws.Range("H96").Value = "17771.00"
$ws.Range("H97").Value = "1288.84"
$ws.Range("H98").Value = "18562500.00"
$ws.Range("H99").Value = "6500.00"
$ws.Range("H100").Value = "35000.00"
$ws.Range("H101").Value = "33000.00"
$ws.Range("H102").Value = "22000.00"
$ws.Range("H103").Value = "22000.00"
$ws.Range("H104").Value = "7000.00"
$ws.Range("H105").Value = "24000.00"
$ws.Range("H106").Value = "10000.00"
$ws.Range("H107").Value = "10000.00"
$ws.Range("H108").Value = "12000.00"
$ws.Range("H109").Value = "10000.00"
$ws.Range("H110").Value = "10000.00"
$ws.Range("H111").Value = "9000.00"
$ws.Range("H112").Value = "10500.00"
$ws.Range("H113").Value = "4000.00"
$ws.Range("H114").Value = "10000.00"
$ws.Range("H115").Value = "13000.00"
$ws.Range("H116").Value = "23500.00"
$ws.Range("H117").Value = "10000.00"
$ws.Range("H118").Value = "5000.00"
$ws.Range("H119").Value = "45305.00"
$ws.Range("H120").Value = "18000.00"
$ws.Range("H121").Value = "10000.00"
$ws.Range("H122").Value = "12000.00"
$ws.Range("H123").Value = "45000.00"
$ws.Range("H124").Value = "18000.00"
$ws.Range("H125").Value = "10000.00"
$ws.Range("H126").Value = "10000.00"
$ws.Range("H127").Value = "50000.00"
$ws.Range("H128").Value = "45000.00"
$ws.Range("H129").Value = "86800.00"
$ws.Range("H130").Value = "20100.00"
$ws.Range("H131").Value = "53000.00"
$ws.Range("H132").Value = "258893.60"
$ws.Range("H133").Value = "18416.00"
$ws.Range("H134").Value = "15400.00"
$ws.Range("H135").Value = "33150.00"
$ws.Range("H136").Value = "15200.00"
$ws.Range("H137").Value = "13200.00"
$ws.Range("H138").Value = "42365.00"
$ws.Range("H139").Value = "338.54"
$ws.Range("H140").Value = "9370.00"
$ws.Range("H141").Value = "4800.00"
$ws.Range("H142").Value = "4500.00"
$ws.Range("H143").Value = "8500.00"
$ws.Range("H144").Value = "20418.79"
$ws.Range("H145").Value = "1297.64"
$ws.Range("H146").Value = "18810.66"
$ws.Range("H147").Value = "118381.88"
$ws.Range("H148").Value = "9290.00"
$ws.Range("H149").Value = "5200.00"
$ws.Range("H150").Value = "1050.00"
$ws.Range("H151").Value = "23000.00"
$ws.Range("H152").Value = "1825.00"
$ws.Range("H153").Value = "85.00"
$ws.Range("H154").Value = "9825.00"
$ws.Range("H155").Value = "600.00"
$ws.Range("H156").Value = "1527733.44"
$ws.Range("H157").Value = "80000.00"
$ws.Range("H158").Value = "40000.00"
$ws.Range("H159").Value = "40000.00"
$ws.Range("H160").Value = "40000.00"
$ws.Range("H161").Value = "80000.00"
$ws.Range("H162").Value = "40000.00"
$ws.Range("H163").Value = "55000.00"
$ws.Range("H164").Value = "40000.00"
$ws.Range("H165").Value = "40000.00"
$ws.Range("H166").Value = "80000.00"
$ws.Range("H167").Value = "80000.00"
$ws.Range("H168").Value = "7500.00"
$ws.Range("H169").Value = "7850.00"
$ws.Range("H170").Value = "69441.51"
$ws.Range("H171").Value = "8230759.63"
$ws.Range("H172").Value = "800.00"
$ws.Range("H173").Value = "298100.00"
$ws.Range("H174").Value = "312760.00"
$ws.Range("H175").Value = "298100.00"
$ws.Range("H176").Value = "298100.00"
$ws.Range("H177").Value = "301460.00"
$ws.Range("H178").Value = "298100.00"
$ws.Range("H179").Value = "566600.00"
$ws.Range("H180").Value = "298100.00"
$ws.Range("H181").Value = "761490.00"
$ws.Range("H182").Value = "1017000.00"
$ws.Range("H183").Value = "434180.00"
$ws.Range("H184").Value = "298100.00"
$ws.Range("H185").Value = "298100.00"
$ws.Range("H186").Value = "596200.00"
$ws.Range("H187").Value = "550900.00"
$ws.Range("H188").Value = "592360.00"
$ws.Range("H189").Value = "878780.00"
$ws.Range("H190").Value = "566600.00"
$ws.Range("H191").Value = "857980.00"
$ws.Range("H192").Value = "596200.00"
$ws.Range("H193").Value = "327403.33"
$ws.Range("H194").Value = "40680.00"
$ws.Range("H195").Value = "70950.00"
$ws.Range("H196").Value = "2312.34"
$ws.Range("H197").Value = "15200.00"
$ws.Range("H198").Value = "3000.00"
$ws.Range("H199").Value = "315000.00"
$ws.Range("H200").Value = "19331.21"
$ws.Range("H201").Value = "20388.02"
$ws.Range("H202").Value = "29000.00"
$ws.Range("H203").Value = "1535.74"
$ws.Range("H204").Value = "1700.00"
$ws.Range("H205").Value = "5542.00"
$ws.Range("H206").Value = "6655.00"
$ws.Range("H207").Value = "1098000.00"
$ws.Range("H208").Value = "1570000.00"
$ws.Range("H209").Value = "945000.00"

# Restore the default (Normal) cell style now that the text values are set,
# so the cells end up visually identical to how they started.
$importeRange.Style = "Normal"

Write-Host "Done updating shared strings."
